# Generate Report for Handback
# Populates the "ce87676b-a28c-462f-870a-f9d9854f71aa" handback row (row 7)
# on both the zh-cn and de-de status sheets now that a (stale) handback
# file has come back for that item.

$wb = $excel.ActiveWorkbook

# ---- zh-cn sheet --------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Latest Target File (I7) becomes a hyperlink to the handback .md, same as
# the other rows in this column.
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/77c65156946417a722dcf686063d76fb363ac5c4/e2e/ce87676b-a28c-462f-870a-f9d9854f71aa.md",
    "",
    "",
    "ce87676b-a28c-462f-870a-f9d9854f71aa.md"
) | Out-Null

# Latest Handback File / DateTime now have real values.
$wsZh.Range("J7").Value = "ce87676b-a28c-462f-870a-f9d9854f71aa.1eb0f393652325458b2f9702a5b93a6c7679c083.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-29 00:54:10"

# Error Detail: the handback is stale relative to the latest handoff.
$wsZh.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77c65156946417a722dcf686063d76fb363ac5c4/e2e/ce87676b-a28c-462f-870a-f9d9854f71aa.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7de1a938cae6e4eb644e33565c76e0e2604cd20/e2e/ce87676b-a28c-462f-870a-f9d9854f71aa.md."

# ---- de-de sheet ---------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/77c65156946417a722dcf686063d76fb363ac5c4/e2e/ce87676b-a28c-462f-870a-f9d9854f71aa.md",
    "",
    "",
    "ce87676b-a28c-462f-870a-f9d9854f71aa.md"
) | Out-Null

$wsDe.Range("J7").Value = "ce87676b-a28c-462f-870a-f9d9854f71aa.1eb0f393652325458b2f9702a5b93a6c7679c083.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-29 00:54:17"

$wsDe.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77c65156946417a722dcf686063d76fb363ac5c4/e2e/ce87676b-a28c-462f-870a-f9d9854f71aa.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7de1a938cae6e4eb644e33565c76e0e2604cd20/e2e/ce87676b-a28c-462f-870a-f9d9854f71aa.md."

Write-Output "Handback report generated for ce87676b-a28c-462f-870a-f9d9854f71aa (zh-cn, de-de)"
